$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '97.545.53'
$ws.Range("E2").Value = '  -1.26%  '

# Row 3
$ws.Range("D3").Value = '3.393.42'
$ws.Range("E3").Value = '  +2.32%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = "'253.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '

# Row 6
$ws.Range("D6").Value = "'650.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.90%  '

# Row 7
$ws.Range("E7").Value = '  +1.58%  '

# Row 8
$ws.Range("D8").Value = "'0.429"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.59%  '

# Row 9
$ws.Range("E9").Value = '  +6.88%  '

# Row 10
$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.01%  '

# Row 11
$ws.Range("D11").Value = '3.386.34'
$ws.Range("E11").Value = '  +2.17%  '

# Row 12
$ws.Range("D12").Value = "'0.211"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.45%  '

# Row 13
$ws.Range("D13").Value = "'41.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.16%  '

# Row 14
$ws.Range("D14").Value = "'6.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +17.24%  '

# Row 15
$ws.Range("E15").Value = '  +2.86%  '

# Row 16
$ws.Range("D16").Value = '97.154.54'
$ws.Range("E16").Value = '  -1.62%  '

# Row 17
$ws.Range("D17").Value = '4.020.56'
$ws.Range("E17").Value = '  +1.64%  '

# Row 18
$ws.Range("D18").Value = "'8.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +27.46%  '

# Row 19
$ws.Range("D19").Value = '3.399.73'
$ws.Range("E19").Value = '  +2.45%  '

# Row 20
$ws.Range("D20").Value = "'17.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.70%  '

# Row 21
$ws.Range("D21").Value = "'0.499"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +47.14%  '

# Row 22
$ws.Range("D22").Value = "'10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +11.19%  '

# Row 23
$ws.Range("E23").Value = '  -2.22%  '

# Row 24
$ws.Range("D24").Value = "'505.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.62%  '

# Row 25
$ws.Range("E25").Value = '  +0.27%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = "'98.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.26%  '

# Row 27
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").Value = "'6.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '

# Row 28
$ws.Range("D28").Value = "'12.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.74%  '

# Row 29
$ws.Range("D29").Value = '3.577.06'
$ws.Range("E29").Value = '  +2.41%  '

# Row 30
$ws.Range("E30").Value = '  +3.25%  '

# Row 31
$ws.Range("E31").Value = '  +6.66%  '

# Row 32
$ws.Range("D32").Value = "'11.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.27%  '

# Row 33
$ws.Range("D33").Value = "'0.995"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.43%  '

# Row 34
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '

# Row 35
$ws.Range("D35").Value = "'0.563"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +15.40%  '

# Row 36
$ws.Range("D36").Value = "'29.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.63%  '

# Row 37
$ws.Range("E37").Value = '  +13.91%  '

# Row 38
$ws.Range("D38").Value = "'7.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.17%  '

# Row 39
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = "'1.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.55%  '

# Row 40
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = "'522.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.73%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = "'0.152"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.14%  '

# Row 42
$ws.Range("D42").Value = "'24.72"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.11%  '

# Row 43
$ws.Range("D43").Value = "'0.853"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.00%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = "'0.0419"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +17.39%  '

# Row 45
$ws.Range("B45").Value = 'MantraDAO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D45").Value = "'3.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.33%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = "'3.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.78%  '

# Row 47
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = "'5.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.46%  '

# Row 48
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.03%  '

# Row 49
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'8.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.42%  '

# Row 50
$ws.Range("D50").Value = "'1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.33%  '

# Row 51
$ws.Range("E51").Value = '  +2.54%  '
